$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace ANILSON (004385806) with THAIS (005395948), balance 250000
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "005395948"
$ws.Range("B2").Value = "THAIS"
$ws.Range("C2").Value = 250000

# Update Saldo values for the next rows
$ws.Range("C3").Value = 130877.29
$ws.Range("C4").Value = 92137.62
$ws.Range("C5").Value = 79573.4
$ws.Range("C6").Value = 23075.45
